$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column to hold a surrogate-key "ID" field; this
# shifts description/amount/date from A/B/C to B/C/D.
$ws.Columns("A:A").Insert()

# Populate the new ID column (upsert primary key) as plain numbers.
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# The ID column is an internal key - keep it narrow and hidden.
$ws.Columns("A:A").ColumnWidth = 7.14
$ws.Columns("A:A").Hidden = $true

# Upsert: the Groceries row's date moved from 2025-07-02 to 2025-07-03.
$ws.Range("D3").Value = 45841
